# Aggiornamento fino a 20/09/2021
# Append new daily rows (375-385) to Sheet1, continuing the existing series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 374

# New data: date serial, nuovi pos. (B), somma mobile 7gg. (C), somma mobile 7gg. per 100mila abitanti (D)
$newRows = @(
    @(44449, 0, 1, 26.76659528907923),
    @(44450, 1, 1, 26.76659528907923),
    @(44451, 0, 1, 26.76659528907923),
    @(44452, 0, 1, 26.76659528907923),
    @(44453, 0, 1, 26.76659528907923),
    @(44454, 0, 1, 26.76659528907923),
    @(44455, 0, 1, 26.76659528907923),
    @(44456, 0, 1, 26.76659528907923),
    @(44457, 0, 0, 0),
    @(44458, 0, 0, 0),
    @(44459, 0, 0, 0)
)

$srcRow = $ws.Range("A" + $lastRow + ":D" + $lastRow)

$r = $lastRow
foreach ($row in $newRows) {
    $r = $r + 1

    $destRow = $ws.Range("A" + $r + ":D" + $r)
    $srcRow.Copy($destRow)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
